# Fruta / hortaliza, semanal
# Insert a new weekly price record as the new row 117 on the single
# worksheet, pushing the existing rows 117:164 down to 118:165.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 117:164 down one row (new empty row 117, dimension grows to R165).
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new data record.
$ws.Range("A117").Value = 5
$ws.Range("B117").Value = "Macroferia Regional de Talca"
$ws.Range("C117").Value = "Maule"
$ws.Range("D117").Value = 44726
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 100112031
$ws.Range("G117").Value = "Poroto verde"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("M117").Value = 25000
$ws.Range("N117").Value = "$/malla 25 kilos"
$ws.Range("O117").Value = "Región de Arica y Parinacota"
$ws.Range("P117").Value = 1000
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"
